# Daily_Status.xlsx — "Dec 3rd - Status" update
# Appends the Dec-3 status entry (and the following days through Dec 7)
# to the bottom of the daily-status log on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 149 already holds the 03/12/2021 date label (existing last row);
# fill in its DONE / PROGRESS columns with the day's status.
$ws.Range("B149").Value = "Revised basic C programming"
$ws.Range("C149").Value = "LDD-DS-OS concepts"

# Continuation line for 03/12/2021.
$ws.Range("B150").Value = "Interview questions"
$ws.Range("C150").Value = "Interview preparation"

# 04/12/2021 — holiday.
$ws.Range("A151").Value = "04/12/2021"
$ws.Range("B151").Value = "HOLIDAY"

# 05/12/2021 — holiday.
$ws.Range("A152").Value = "05/12/2021"
$ws.Range("B152").Value = "HOLIDAY"

# 06/12/2021 — leave.
$ws.Range("A153").Value = "06/12/2021"
$ws.Range("B153").Value = "LEAVE"

# 07/12/2021 — new entry, status to be filled in later.
$ws.Range("A154").Value = "07/12/2021"

# Leave the selection on the new last row, matching where Excel's cursor
# would land after typing these entries in.
$ws.Range("A154").Select()
